# Updates the Betfair odds/liquidity grid (columns F..AO, rows 2..16 of
# Sheet1) so the numeric values match the refreshed odds captured in the
# commit's XLSX diff. Every entry below is (row, column, new value); the
# column is given as its 1-based index (F=6 ... AO=41).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{R=2; C=6; V=2.8},
    @{R=2; C=7; V=2.9},
    @{R=2; C=10; V=3.6},
    @{R=2; C=11; V=3.8},
    @{R=2; C=14; V=4.4},
    @{R=2; C=16; V=2.18},
    @{R=2; C=17; V=1.79},
    @{R=2; C=18; V=1.47},
    @{R=2; C=19; V=2.96},
    @{R=2; C=20; V=1.65},
    @{R=2; C=22; V=1.6},
    @{R=2; C=24; V=20},
    @{R=2; C=25; V=13.5},
    @{R=2; C=26; V=22},
    @{R=2; C=27; V=1000},
    @{R=2; C=28; V=14},
    @{R=2; C=29; V=8.4},
    @{R=2; C=30; V=13},
    @{R=2; C=31; V=28},
    @{R=2; C=32; V=23},
    @{R=2; C=33; V=13.5},
    @{R=2; C=34; V=16.5},
    @{R=2; C=35; V=36},
    @{R=2; C=38; V=1000},
    @{R=2; C=39; V=70},
    @{R=2; C=40; V=23},
    @{R=2; C=41; V=19.5},
    @{R=3; C=6; V=13},
    @{R=3; C=7; V=28},
    @{R=3; C=9; V=1.21},
    @{R=3; C=11; V=16.5},
    @{R=3; C=12; V=1.17},
    @{R=3; C=14; V=7.6},
    @{R=3; C=15; V=1.11},
    @{R=3; C=16; V=3.5},
    @{R=3; C=17; V=1.29},
    @{R=3; C=18; V=2.02},
    @{R=3; C=19; V=1.72},
    @{R=3; C=20; V=2.06},
    @{R=3; C=21; V=1.78},
    @{R=3; C=23; V=1.04},
    @{R=3; C=24; V=60},
    @{R=3; C=25; V=970},
    @{R=3; C=26; V=970},
    @{R=3; C=27; V=970},
    @{R=3; C=28; V=90},
    @{R=3; C=29; V=27},
    @{R=3; C=30; V=970},
    @{R=3; C=31; V=970},
    @{R=3; C=33; V=85},
    @{R=3; C=34; V=48},
    @{R=3; C=35; V=48},
    @{R=3; C=41; V=2.98},
    @{R=4; C=6; V=2.56},
    @{R=4; C=7; V=2.96},
    @{R=4; C=8; V=2.52},
    @{R=4; C=9; V=2.88},
    @{R=4; C=10; V=3.5},
    @{R=4; C=11; V=4.1},
    @{R=4; C=12; V=1.33},
    @{R=4; C=13; V=1.05},
    @{R=4; C=14; V=4.2},
    @{R=4; C=15; V=1.24},
    @{R=4; C=16; V=2.12},
    @{R=4; C=17; V=1.72},
    @{R=4; C=18; V=1.44},
    @{R=4; C=19; V=2.8},
    @{R=4; C=20; V=1.61},
    @{R=4; C=21; V=2.28},
    @{R=4; C=22; V=1.53},
    @{R=4; C=23; V=1.51},
    @{R=4; C=24; V=23},
    @{R=4; C=25; V=970},
    @{R=4; C=26; V=21},
    @{R=4; C=27; V=42},
    @{R=4; C=28; V=970},
    @{R=4; C=29; V=970},
    @{R=4; C=30; V=970},
    @{R=4; C=31; V=29},
    @{R=4; C=32; V=21},
    @{R=4; C=33; V=970},
    @{R=4; C=34; V=970},
    @{R=4; C=35; V=38},
    @{R=4; C=36; V=42},
    @{R=4; C=37; V=30},
    @{R=4; C=38; V=38},
    @{R=4; C=39; V=85},
    @{R=4; C=40; V=21},
    @{R=4; C=41; V=21},
    @{R=5; C=6; V=3.85},
    @{R=5; C=7; V=5.2},
    @{R=5; C=8; V=1.81},
    @{R=5; C=9; V=2},
    @{R=5; C=10; V=3.9},
    @{R=5; C=12; V=1.29},
    @{R=5; C=15; V=1.2},
    @{R=5; C=16; V=2.36},
    @{R=5; C=19; V=2.46},
    @{R=5; C=20; V=1.59},
    @{R=5; C=21; V=2.34},
    @{R=5; C=22; V=2},
    @{R=5; C=23; V=1.27},
    @{R=5; C=26; V=17},
    @{R=5; C=27; V=26},
    @{R=5; C=32; V=42},
    @{R=5; C=33; V=21},
    @{R=5; C=34; V=20},
    @{R=5; C=36; V=100},
    @{R=5; C=37; V=55},
    @{R=5; C=39; V=80},
    @{R=5; C=40; V=42},
    @{R=6; C=7; V=3.7},
    @{R=6; C=17; V=1.95},
    @{R=6; C=20; V=1.77},
    @{R=6; C=21; V=2.2},
    @{R=7; C=6; V=2.44},
    @{R=7; C=7; V=2.74},
    @{R=7; C=8; V=2.84},
    @{R=7; C=9; V=3.25},
    @{R=7; C=11; V=3.85},
    @{R=7; C=13; V=1.07},
    @{R=7; C=17; V=1.89},
    @{R=7; C=21; V=2.14},
    @{R=7; C=25; V=14.5},
    @{R=7; C=26; V=24},
    @{R=7; C=27; V=60},
    @{R=7; C=29; V=8.4},
    @{R=7; C=30; V=15.5},
    @{R=7; C=31; V=40},
    @{R=7; C=34; V=19.5},
    @{R=7; C=39; V=100},
    @{R=7; C=41; V=32},
    @{R=8; C=6; V=2.26},
    @{R=8; C=7; V=2.7},
    @{R=8; C=8; V=2.98},
    @{R=8; C=9; V=3.75},
    @{R=8; C=10; V=3.35},
    @{R=8; C=11; V=4},
    @{R=8; C=13; V=1.06},
    @{R=8; C=17; V=1.71},
    @{R=8; C=20; V=1.67},
    @{R=8; C=21; V=2.16},
    @{R=8; C=24; V=19.5},
    @{R=8; C=25; V=16.5},
    @{R=8; C=26; V=28},
    @{R=8; C=27; V=65},
    @{R=8; C=28; V=13.5},
    @{R=8; C=29; V=10},
    @{R=8; C=30; V=17},
    @{R=8; C=31; V=44},
    @{R=8; C=32; V=970},
    @{R=8; C=33; V=14},
    @{R=8; C=34; V=20},
    @{R=8; C=35; V=55},
    @{R=8; C=36; V=40},
    @{R=8; C=37; V=32},
    @{R=8; C=38; V=44},
    @{R=8; C=39; V=100},
    @{R=8; C=40; V=23},
    @{R=8; C=41; V=38},
    @{R=9; C=6; V=1.78},
    @{R=9; C=7; V=1.86},
    @{R=9; C=8; V=4.6},
    @{R=9; C=10; V=3.9},
    @{R=9; C=17; V=1.81},
    @{R=9; C=20; V=1.78},
    @{R=9; C=24; V=20},
    @{R=9; C=28; V=11},
    @{R=10; C=8; V=2.36},
    @{R=10; C=9; V=2.58},
    @{R=10; C=10; V=3.4},
    @{R=10; C=11; V=3.8},
    @{R=10; C=13; V=1.06},
    @{R=10; C=16; V=1.96},
    @{R=10; C=17; V=1.74},
    @{R=10; C=20; V=1.7},
    @{R=10; C=21; V=2.2},
    @{R=10; C=24; V=970},
    @{R=10; C=25; V=13.5},
    @{R=10; C=26; V=20},
    @{R=10; C=28; V=16},
    @{R=10; C=29; V=9.8},
    @{R=10; C=32; V=26},
    @{R=10; C=33; V=16.5},
    @{R=10; C=35; V=42},
    @{R=10; C=37; V=40},
    @{R=11; C=6; V=1.63},
    @{R=11; C=8; V=5},
    @{R=11; C=9; V=5.8},
    @{R=11; C=11; V=4.8},
    @{R=11; C=13; V=1.04},
    @{R=11; C=17; V=1.58},
    @{R=11; C=20; V=1.71},
    @{R=11; C=24; V=22},
    @{R=11; C=25; V=26},
    @{R=11; C=27; V=140},
    @{R=11; C=33; V=11.5},
    @{R=12; C=10; V=3.8},
    @{R=12; C=21; V=2.46},
    @{R=12; C=25; V=18},
    @{R=12; C=29; V=8.8},
    @{R=12; C=33; V=11},
    @{R=12; C=37; V=20},
    @{R=13; C=6; V=5.9},
    @{R=13; C=7; V=6.2},
    @{R=13; C=8; V=1.6},
    @{R=13; C=9; V=1.63},
    @{R=13; C=10; V=4.5},
    @{R=13; C=11; V=4.8},
    @{R=13; C=16; V=2.64},
    @{R=13; C=17; V=1.57},
    @{R=13; C=20; V=1.66},
    @{R=13; C=21; V=2.4},
    @{R=13; C=24; V=27},
    @{R=13; C=25; V=12},
    @{R=13; C=26; V=12.5},
    @{R=13; C=27; V=17},
    @{R=13; C=28; V=29},
    @{R=13; C=29; V=11},
    @{R=13; C=31; V=16},
    @{R=13; C=33; V=24},
    @{R=13; C=34; V=18.5},
    @{R=13; C=35; V=27},
    @{R=13; C=36; V=150},
    @{R=13; C=37; V=70},
    @{R=13; C=38; V=60},
    @{R=13; C=39; V=75},
    @{R=13; C=40; V=60},
    @{R=13; C=41; V=6.4},
    @{R=14; C=6; V=3.8},
    @{R=14; C=7; V=3.9},
    @{R=14; C=8; V=2},
    @{R=14; C=9; V=2.04},
    @{R=14; C=10; V=4},
    @{R=14; C=16; V=2.48},
    @{R=14; C=21; V=2.6},
    @{R=14; C=26; V=15},
    @{R=14; C=29; V=9.4},
    @{R=14; C=31; V=19},
    @{R=14; C=34; V=16},
    @{R=14; C=35; V=29},
    @{R=14; C=41; V=9.6},
    @{R=15; C=6; V=1.91},
    @{R=15; C=8; V=4.6},
    @{R=15; C=16; V=1.59},
    @{R=15; C=17; V=2.44},
    @{R=16; C=6; V=3.65},
    @{R=16; C=7; V=4.4},
    @{R=16; C=8; V=2.14},
    @{R=16; C=9; V=2.46},
    @{R=16; C=10; V=3},
    @{R=16; C=11; V=3.5},
    @{R=16; C=13; V=1.09},
    @{R=16; C=16; V=1.57},
    @{R=16; C=17; V=2.4},
    @{R=16; C=20; V=2.02},
    @{R=16; C=21; V=1.79},
    @{R=16; C=24; V=11.5},
    @{R=16; C=25; V=9},
    @{R=16; C=26; V=970},
    @{R=16; C=27; V=40},
    @{R=16; C=28; V=970},
    @{R=16; C=29; V=9},
    @{R=16; C=30; V=970},
    @{R=16; C=31; V=38},
    @{R=16; C=32; V=34},
    @{R=16; C=33; V=21},
    @{R=16; C=34; V=28},
    @{R=16; C=36; V=110},
    @{R=16; C=38; V=100},
    @{R=16; C=39; V=210}
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.R, $u.C).Value = $u.V
}
